$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-31 16:52:40"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-31 16:52:35"
$wsZhCn.Range("K4").Value = "2016-08-31 16:53:11"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-31 16:53:19"
